# ---------------------------------------------------------------------------
# Commit: "feat: add 2022-Q4 data"
#
# 1) On the "总计" (summary) sheet, insert a new row 2 for the 2022-Q4
#    quarter (pushing every existing quarter row down by one) and
#    renumber the running index in column A.
# 2) Insert a brand-new worksheet named "2022-Q4" right after "总计"
#    (and before "2022-Q3") holding the per-fund breakdown for that
#    quarter.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ===========================================================================
# Part 1 - "总计" summary sheet: insert the 2022-Q4 row
# ===========================================================================
$summary = $wb.Worksheets.Item(1)

# Push rows 2..9 down to 3..10, duplicating the formatting of the old row 2
# (now row 3) onto the freshly inserted row 2.
$summary.Rows.Item(2).Insert()
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value2 = 0
$summary.Cells.Item(2, 2).Value2 = "2022-Q4"
$summary.Cells.Item(2, 3).Value2 = 9
$summary.Cells.Item(2, 4).Value2 = 5.9

# Renumber the running index (column A) for the rows that shifted down.
for ($r = 3; $r -le 10; $r++) {
    $summary.Cells.Item($r, 1).Value2 = $r - 2
}

# ===========================================================================
# Part 2 - new "2022-Q4" detail sheet
# ===========================================================================
# The existing "2020-Q4" sheet (last tab) already has the exact column
# layout/styling we need (基金代码/基金名称/.../仓位排名), so clone it into
# the right tab position, rename it, trim it to 9 data rows, then overwrite
# every cell with the 2022-Q4 figures.
$template = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $summary)

$detail = $wb.Worksheets.Item(2)
$detail.Name = "2022-Q4"

# Only 9 funds this quarter (old template sheet had 13 data rows).
$detail.Range("A11:H14").Delete(-4162)

# Header row (column D label differs: 基金规模 instead of 基金金额).
$detail.Cells.Item(1, 4).Formula = "'基金规模"

$f2  = @("009076", "工银圆兴混合",                       "49.11", "92.89", "3.38", "1.6599", 8)
$f3  = @("010591", "富国中国中小盘混合（QDII）美元",        "32.21", "87.21", "4.65", "1.4978", 2)
$f4  = @("100061", "富国中国中小盘混合（QDII）人民币",       "32.21", "87.21", "4.65", "1.4978", 2)
$f5  = @("009029", "工银高质量成长混合A",                  "12.82", "87.91", "3.46", "0.4436", 7)
$f6  = @("006752", "天弘港股通精选灵活配置混合A",           "5.05",  "89.88", "7.64", "0.3858", 4)
$f7  = @("006753", "天弘港股通精选灵活配置混合C",           "2.97",  "89.88", "7.64", "0.2269", 4)
$f8  = @("007109", "南方沪港深核心优势混合",                "1.99",  "89.23", "3.30", "0.0657", 10)
$f9  = @("160125", "南方香港优选股票（QDII-LOF）",          "2.27",  "84.75", "2.72", "0.0617", 10)
$f10 = @("009030", "工银高质量成长混合C",                  "1.72",  "87.91", "3.46", "0.0595", 7)

$funds = @($f2, $f3, $f4, $f5, $f6, $f7, $f8, $f9, $f10)

$r = 2
foreach ($fund in $funds) {
    $detail.Cells.Item($r, 1).Value2 = $r - 2
    $detail.Cells.Item($r, 2).Formula = "'" + $fund[0]
    $detail.Cells.Item($r, 3).Formula = "'" + $fund[1]
    $detail.Cells.Item($r, 4).Formula = "'" + $fund[2]
    $detail.Cells.Item($r, 5).Formula = "'" + $fund[3]
    $detail.Cells.Item($r, 6).Formula = "'" + $fund[4]
    $detail.Cells.Item($r, 7).Formula = "'" + $fund[5]
    $detail.Cells.Item($r, 8).Value2 = $fund[6]
    $r = $r + 1
}
